$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; B = 4807; C = 3975; D = 1.54785943031311 },
    @{ Row = 3; B = 4723; C = 3513; D = 1.54785943031311 },
    @{ Row = 4; B = 5270; C = 5279; D = 1.54785943031311 },
    @{ Row = 5; B = 4575; C = 5171; D = 1.54785943031311 },
    @{ Row = 6; B = 5822; C = 4914; D = 1.54785943031311 },
    @{ Row = 7; B = 9789; C = 7223; D = 1.54785943031311 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.D
}
